$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Metadata sheet value updates
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Metadata")

# Experimental (B7): false -> true
# Write via a text formula then paste-special as values so the result lands
# as a genuine shared string ("true") instead of Excel's auto-detected
# Boolean literal.
$b7 = $ws.Range("B7")
$b7.Formula = "=""true"""
$b7.Copy()
$b7.PasteSpecial(-4163)

# Date (B8): refreshed timestamp
$ws.Range("B8").Value = "2024-02-19T18:37:26-06:00"

# Case Sensitive (B14): (empty) -> true
$b14 = $ws.Range("B14")
$b14.Formula = "=""true"""
$b14.Copy()
$b14.PasteSpecial(-4163)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Re-apply (and mark as explicitly applied) the vertical-top / wrap-text
# alignment used by the header and body styles on both sheets.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Concepts")

$ranges = @(
    $ws.Range("A1:B1"),
    $ws.Range("A2:B21"),
    $ws2.Range("A1:D1"),
    $ws2.Range("A2:D6")
)

foreach ($rng in $ranges) {
    $rng.VerticalAlignment = -4160
    $rng.WrapText = $true
}
